$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for Frambuesa (raspberry) entries get cyclically rotated:
# new row3 <- old row6, new row4 <- old row3, new row6 <- old row4
# (columns D, M, N, O, P, R, S)

$oldRow3 = @{
    D = $ws.Range("D3").Value2()
    M = $ws.Range("M3").Value2()
    N = $ws.Range("N3").Value2()
    O = $ws.Range("O3").Value2()
    P = $ws.Range("P3").Value2()
    R = $ws.Range("R3").Value2()
    S = $ws.Range("S3").Value2()
}

$oldRow4 = @{
    D = $ws.Range("D4").Value2()
    M = $ws.Range("M4").Value2()
    N = $ws.Range("N4").Value2()
    O = $ws.Range("O4").Value2()
    P = $ws.Range("P4").Value2()
    R = $ws.Range("R4").Value2()
    S = $ws.Range("S4").Value2()
}

$oldRow6 = @{
    D = $ws.Range("D6").Value2()
    M = $ws.Range("M6").Value2()
    N = $ws.Range("N6").Value2()
    O = $ws.Range("O6").Value2()
    P = $ws.Range("P6").Value2()
    R = $ws.Range("R6").Value2()
    S = $ws.Range("S6").Value2()
}

# Row 3 <- old row 6
$ws.Range("D3").Value = $oldRow6.D
$ws.Range("M3").Value = $oldRow6.M
$ws.Range("N3").Value = $oldRow6.N
$ws.Range("O3").Value = $oldRow6.O
$ws.Range("P3").Value = $oldRow6.P
$ws.Range("R3").Value = $oldRow6.R
$ws.Range("S3").Value = $oldRow6.S

# Row 4 <- old row 3
$ws.Range("D4").Value = $oldRow3.D
$ws.Range("M4").Value = $oldRow3.M
$ws.Range("N4").Value = $oldRow3.N
$ws.Range("O4").Value = $oldRow3.O
$ws.Range("P4").Value = $oldRow3.P
$ws.Range("R4").Value = $oldRow3.R
$ws.Range("S4").Value = $oldRow3.S

# Row 6 <- old row 4
$ws.Range("D6").Value = $oldRow4.D
$ws.Range("M6").Value = $oldRow4.M
$ws.Range("N6").Value = $oldRow4.N
$ws.Range("O6").Value = $oldRow4.O
$ws.Range("P6").Value = $oldRow4.P
$ws.Range("R6").Value = $oldRow4.R
$ws.Range("S6").Value = $oldRow4.S
